$d = $word.ActiveDocument

# 1. Remove the "Meta description: ..." paragraph that currently sits right
#    after the title heading (it is being moved/reworked further down).
$metaPara = $d.Paragraphs.Item(2)
[void]$metaPara.Range.Delete()

# 2. Insert a new bold paragraph ("Play 15 Armadillos Slot Free - Features
#    Exciting Bonus Rounds") right before the final ("Prompt: ...") paragraph.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
[void]$lastPara.Range.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item($count)
[void]$newPara.Range.InsertXML("<w:p xmlns:w=""http://schemas.openxmlformats.org/wordprocessingml/2006/main""><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play 15 Armadillos Slot Free - Features Exciting Bonus Rounds</w:t></w:r></w:p>")

# 3. Replace the text of the final paragraph (previously the "Prompt: ..."
#    text) with the former meta description text, keeping its italic run
#    formatting intact.
$count = $d.Paragraphs.Count
$finalPara = $d.Paragraphs.Item($count)
$oldText = 'Prompt: Create a cartoon-style feature image for the game "15 Armadillos" that features a happy Maya warrior with glasses. For the feature image of "15 Armadillos", let''s have a cartoon-style design featuring a happy Maya warrior with glasses. The warrior can be depicted wearing a headdress made of colorful feathers, with intricate designs on their face and body. They can be holding a staff or weapon made of stone or wood, with a happy expression on their face. In the background, we can see the Everglades National Park with its lush greenery and animals like alligators and otters. The image can be bright and colorful to reflect the fun and adventurous nature of the game.'
$newText = "Read our review and play 15 Armadillos slot for free. Enjoy exciting bonus rounds such as Wild Respins, Armadillo Link, and Free Spins."
[void]$finalPara.Range.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
